$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '68.950.33'
$ws.Range('E2').Value = '  -2.75%  '

Set-TextValue 'D3' '3.677.55'
$ws.Range('E3').Value = '  -3.75%  '

Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.07%  '

Set-TextValue 'D5' '678.86'
$ws.Range('E5').Value = '  -3.98%  '

Set-TextValue 'D6' '162.03'
$ws.Range('E6').Value = '  -4.83%  '

Set-TextValue 'D7' '3.677.69'
$ws.Range('E7').Value = '  -3.67%  '

Set-TextValue 'D8' '1.00'
$ws.Range('E8').Value = '  -0.08%  '

Set-TextValue 'D9' '0.491'
$ws.Range('E9').Value = '  -6.00%  '

$ws.Range('E10').Value = '  -8.31%  '

$ws.Range('E11').Value = '  -2.10%  '

Set-TextValue 'D12' '0.448'
$ws.Range('E12').Value = '  -1.82%  '

$ws.Range('E13').Value = '  -7.55%  '

Set-TextValue 'D14' '33.29'
$ws.Range('E14').Value = '  -8.49%  '

Set-TextValue 'D15' '4.296.79'
$ws.Range('E15').Value = '  -3.82%  '

Set-TextValue 'D16' '3.669.24'
$ws.Range('E16').Value = '  -6.84%  '

Set-TextValue 'D17' '69.021.47'
$ws.Range('E17').Value = '  -2.75%  '

$ws.Range('E18').Value = '  -1.70%  '

Set-TextValue 'D19' '16.27'
$ws.Range('E19').Value = '  -5.92%  '

Set-TextValue 'D20' '6.56'
$ws.Range('E20').Value = '  -8.66%  '

Set-TextValue 'D21' '479.83'
$ws.Range('E21').Value = '  -2.80%  '

Set-TextValue 'D22' '9.80'
$ws.Range('E22').Value = '  -7.53%  '

Set-TextValue 'D23' '0.662'
$ws.Range('E23').Value = '  -9.37%  '

Set-TextValue 'D24' '79.09'
$ws.Range('E24').Value = '  -7.65%  '

Set-TextValue 'D25' '3.815.71'
$ws.Range('E25').Value = '  -4.03%  '

$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D26' '0.0000128'
$ws.Range('E26').Value = '  -11.31%  '

$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D27' '11.54'
$ws.Range('E27').Value = '  -4.36%  '

Set-TextValue 'D28' '0.999'
$ws.Range('E28').Value = '  -0.02%  '

Set-TextValue 'D29' '9.43'
$ws.Range('E29').Value = '  -10.56%  '

$ws.Range('E30').Value = '  -13.10%  '

$ws.Range('E31').Value = '  -11.71%  '

Set-TextValue 'D32' '2.09'
$ws.Range('E32').Value = '  -5.88%  '

Set-TextValue 'D33' '6.70'
$ws.Range('E33').Value = '  -9.17%  '

$ws.Range('E34').Value = '  -0.06%  '

$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D35' '26.63'
$ws.Range('E35').Value = '  -8.98%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D36' '0.163'
$ws.Range('E36').Value = '  -6.22%  '

Set-TextValue 'D37' '3.635.11'
$ws.Range('E37').Value = '  -4.09%  '

$ws.Range('E38').Value = '  -7.02%  '

$ws.Range('E39').Value = '  +1.85%  '

$ws.Range('E40').Value = '  -8.57%  '

$ws.Range('E41').Value = '  -0.02%  '

Set-TextValue 'D42' '2.19'
$ws.Range('E42').Value = '  -5.19%  '

$ws.Range('E43').Value = '  -0.11%  '

Set-TextValue 'D44' '0.954'
$ws.Range('E44').Value = '  -8.97%  '

Set-TextValue 'D45' '160.15'
$ws.Range('E45').Value = '  -2.29%  '

Set-TextValue 'D46' '48.25'
$ws.Range('E46').Value = '  -1.19%  '

$ws.Range('E47').Value = '  -13.35%  '

Set-TextValue 'D48' '1.30'
$ws.Range('E48').Value = '  -4.05%  '

$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D49' '389.18'
$ws.Range('E49').Value = '  -8.78%  '

$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue 'D50' '0.000273'
$ws.Range('E50').Value = '  -11.83%  '

Set-TextValue 'D51' '8.00'
$ws.Range('E51').Value = '  -8.38%  '
